$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7925423383712769
$ws.Range("B1").Value = 1.213666796684265
$ws.Range("C1").Value = 2.465372800827026
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.761898636817932
